# Apply the 2022 Tab06 update:
#  - Fix "Etats" -> "États" typos in several region labels
#  - Insert the missing "RDM, pays en développement sans littoral" label
#    (was previously duplicated with "Afrique, pays en développement sans littoral")
#  - Update the "Responsabilité" disclaimer wording
#  - Refresh LLDC (landlocked developing countries) indicator figures in row 92

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row labels (column B / A) ---------------------------------------
$ws.Range("B93").Value = "Afrique, petits États insulaires en développement"
$ws.Range("B94").Value = "RDM, petits États insulaires en développement"
# B95 ("Afrique, pays en développement sans littoral") is unchanged.
$ws.Range("B96").Value = "RDM, pays en développement sans littoral"
$ws.Range("B97").Value = "Afrique, États fragiles"
$ws.Range("B98").Value = "RDM, États fragiles"

# --- Footnote wording ---------------------------------------------------
$ws.Range("A104").Value = "Responsabilité : Ce tableau ainsi que toutes les données qu'il peut comprendre, sont sans préjudice du statut de tout territoire, de la souveraineté s'exerçant sur ce dernier, du tracé des frontières et limites internationales, et du nom de tout territoire, ville ou région."

# --- Updated indicator values for row 92 (RDM, pays les moins avancés) --
$ws.Range("C92").Value = 85.9510744444445
$ws.Range("D92").Value = 84.18435
$ws.Range("E92").Value = 87.7443144444445
$ws.Range("F92").Value = 0.94400555555556
$ws.Range("G92").Value = 69.4269988888889
$ws.Range("H92").Value = 63.2261833333334
$ws.Range("I92").Value = 75.8753355555556
$ws.Range("J92").Value = 0.81008444444444
